# FileCounter -> Dummy: rename the scanner reported in the recorded scan
# issues ("FileCounter" was the superfluous file-counting scanner that was
# replaced by a "Dummy" scanner used for the new ScannerIntegrationFunTest).
# This shows up in the "Scan Issues" cell of the Summary sheet (aggregated,
# multi-line text prefixed with the package id) and of the package detail
# sheet (single-line text).

$wb = $excel.ActiveWorkbook
$wsSummary = $wb.Worksheets.Item(1)
$wsGradle  = $wb.Worksheets.Item(2)

$summaryIssueText = "Gradle:org.ossreviewtoolkit.gradle.example:lib:1.0.0`n  Unknown time [ERROR]: Dummy - Could not download 'Gradle:org.ossreviewtoolkit.gradle.example:lib:1.0.0': DownloadException: Download failed for 'Gradle:org.ossreviewtoolkit.gradle.example:lib:1.0.0'.`nSuppressed: DownloadException: No VCS URL provided for 'Gradle:org.ossreviewtoolkit.gradle.example:lib:1.0.0'. Please make sure the published POM file includes the SCM connection, see: https://docs.gradle.org/current/userguide/publishing_maven.html#sec:modifying_the_generated_pom`nSuppressed: DownloadException: No source artifact URL provided for 'Gradle:org.ossreviewtoolkit.gradle.example:lib:1.0.0'.`n"
$wsSummary.Range("F11").Value = $summaryIssueText

$gradleIssueText = "Unknown time [ERROR]: Dummy - Could not download 'Gradle:org.ossreviewtoolkit.gradle.example:lib:1.0.0': DownloadException: Download failed for 'Gradle:org.ossreviewtoolkit.gradle.example:lib:1.0.0'.`nSuppressed: DownloadException: No VCS URL provided for 'Gradle:org.ossreviewtoolkit.gradle.example:lib:1.0.0'. Please make sure the published POM file includes the SCM connection, see: https://docs.gradle.org/current/userguide/publishing_maven.html#sec:modifying_the_generated_pom`nSuppressed: DownloadException: No source artifact URL provided for 'Gradle:org.ossreviewtoolkit.gradle.example:lib:1.0.0'."
$wsGradle.Range("F11").Value = $gradleIssueText

# Restore the package detail sheet's own last selection (F11), then switch
# back to - and leave selected on - the Summary sheet, which becomes the
# active tab again (matching the saved workbook/view state).
$wsGradle.Range("F11").Select()
$wsSummary.Range("F20").Select()
